$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.018.86'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '1.827.35'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6188'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.98%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.47'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07487'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2907'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.59'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07612'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('D13').Value = '1.828.78'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6608'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009053'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +8.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.957'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('D19').Value = '29.019.12'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('D20').Value = '2.076.57'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '223.99'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.154'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.001'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.378'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1349'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.51%  '
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.495'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.204'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.023'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.035'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05199'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.825'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7305'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.145'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.646'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').Value = '1.271.43'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01777'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.312'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8916'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').Value = '1.974.94'
$ws.Range('E46').Value = '  -0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5116'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.24'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.680'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3947'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.82%  '
